$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("计算器")
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 55
